$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...then you must add six framework names..." -> "...eight..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("then you must add six", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "then you must add eight", 2)

# ---------------------------------------------------------------------------
# Change 2: the sample line gains two more "|no" fields
#   MY ROM|my|framework-res_my.apk|com.htc.resources_my.apk|no|no|no|no
#     -> ...|no|no|no|no|no|no
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("com.htc.resources_my.apk|no|no|no|no", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "com.htc.resources_my.apk|no|no|no|no|no|no", 2)

# ---------------------------------------------------------------------------
# Change 3: the paragraph break that used to sit right after the "...no"
# list (i.e. right before "-Now open program, ") moves down: the sample
# line's paragraph now swallows the start of the next paragraph up through
# the bookmark, and "-Now open program, " becomes the first sentence of a
# fresh paragraph.
# ---------------------------------------------------------------------------

# 3a) Merge the two paragraphs by deleting the paragraph mark that precedes
#     "-Now open pro..."
$full = $d.Content.Text
$idxPro = $full.IndexOf("-Now open pro")
$pmark = $d.Range($idxPro - 1, $idxPro)
$pmark.Delete()

# 3b) Re-insert a paragraph break in the same textual spot (immediately
#     before "-Now"), which now falls right after the lengthened "no" list.
$full = $d.Content.Text
$noListEnd = $full.IndexOf("no|no|no|no|no|no") + ("no|no|no|no|no|no").Length
$breakPoint = $d.Range($noListEnd, $noListEnd)
$breakPoint.InsertParagraphBefore()

# 3c) The "_GoBack" bookmark now sits one character earlier than before -
#     between the last "n" and the trailing "o" of that same "no" list,
#     rather than inside "program". Re-seat it there (Bookmarks.Add moves
#     an existing bookmark of the same name).
$bmPoint = $d.Range($noListEnd - 1, $noListEnd - 1)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# 3d) Normalize "-Now open pro" + "gram, " back into one run of
#     "-Now open program, " text (content is unchanged, only how it is
#     split across runs) without touching the following "goto" run.
$full = $d.Content.Text
$idxPro = $full.IndexOf("-Now open pro")
$idxAfterGram = $full.IndexOf("gram, ") + ("gram, ").Length
$mergeRange = $d.Range($idxPro, $idxAfterGram)
$mergeRange.Text = "@@TMP_MERGE@@"
$full = $d.Content.Text
$idxTmp = $full.IndexOf("@@TMP_MERGE@@")
$tmpRange = $d.Range($idxTmp, $idxTmp + ("@@TMP_MERGE@@").Length)
$tmpRange.Text = "-Now open program, "
